# Auto-generated script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.101.68"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3
$ws.Range("D3").Value = "3.039.38"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "513.85"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "140.73"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.438"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -4.89%  "

# Row 10
$ws.Range("E10").Value = "  -1.02%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.377"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.81%  "

# Row 12
$ws.Range("D12").Value = "3.569.11"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13
$ws.Range("E13").Value = "  -3.22%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.95"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.19%  "

# Row 15
$ws.Range("E15").Value = "  +1.61%  "

# Row 16
$ws.Range("D16").Value = "57.160.06"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.80%  "

# Row 18
$ws.Range("D18").Value = "3.038.73"
$ws.Range("E18").Value = "  +0.59%  "

# Row 19
$ws.Range("E19").Value = "  +4.34%  "

# Row 20
$ws.Range("E20").Value = "  +2.00%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "329.99"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("E23").Value = "  +1.35%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.38"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.56%  "

# Row 25
$ws.Range("D25").Value = "3.166.10"
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("E27").Value = "  -1.65%  "

# Row 28
$ws.Range("E28").Value = "  -3.76%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.72"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.91%  "

# Row 30
$ws.Range("E30").Value = "  -2.23%  "

# Row 31
$ws.Range("E31").Value = "  -0.08%  "

# Row 32
$ws.Range("E32").Value = "  +1.10%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "20.75"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "

# Row 34
$ws.Range("E34").Value = "  -1.82%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "152.50"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "

# Row 37
$ws.Range("E37").Value = "  -0.28%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "25.21"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.38%  "

# Row 39
$ws.Range("E39").Value = "  -0.20%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.89"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.75%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.71"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.96%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("E43").Value = "  +1.62%  "

# Row 44
$ws.Range("E44").Value = "  -1.38%  "

# Row 45
$ws.Range("D45").Value = "2.198.69"
$ws.Range("E45").Value = "  -1.03%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "6.10"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.67%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.947"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.52%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0242"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "

# Row 49
$ws.Range("E49").Value = "  +2.63%  "

# Row 50
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("B51").Value = "Notcoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/2L2Y4ghjj+notcoin-not"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0171"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.73%  "

Write-Output "Applied crypto list refresh"
